$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.4
$ws.Range("J2").Value = 3.15
$ws.Range("L2").Value = 1.49
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 3.1
$ws.Range("T2").Value = 1.91
$ws.Range("W2").Value = 1.58
$ws.Range("AB2").Value = 8.800000000000001
$ws.Range("AK2").Value = 34
$ws.Range("G3").Value = 1.81
$ws.Range("I3").Value = 7
$ws.Range("L3").Value = 1.54
$ws.Range("P3").Value = 1.64
$ws.Range("Q3").Value = 2.28
$ws.Range("S3").Value = 4.5
$ws.Range("T3").Value = 2.14
$ws.Range("Y3").Value = 18.5
$ws.Range("AB3").Value = 6.8
$ws.Range("AF3").Value = 9.800000000000001
$ws.Range("AN3").Value = 1000
$ws.Range("F4").Value = 2.94
$ws.Range("H4").Value = 2.64
$ws.Range("I4").Value = 2.92
$ws.Range("T4").Value = 1.9
$ws.Range("AB4").Value = 12
$ws.Range("AF4").Value = 24
$ws.Range("AG4").Value = 17
$ws.Range("F5").Value = 1.64
$ws.Range("G5").Value = 1.73
$ws.Range("J5").Value = 3.75
$ws.Range("P5").Value = 1.72
$ws.Range("R5").Value = 1.26
$ws.Range("T5").Value = 2.14
$ws.Range("W5").Value = 2.36
$ws.Range("AG5").Value = 12
$ws.Range("AL5").Value = 60
$ws.Range("AN5").Value = 16
$ws.Range("G6").Value = 1.78
$ws.Range("H6").Value = 2.32
$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 3.4
$ws.Range("K6").Value = 980
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 1.26
$ws.Range("P6").Value = 1.24
$ws.Range("Q6").Value = 1.36
$ws.Range("R6").Value = 1.18
$ws.Range("T6").Value = 1.03
$ws.Range("U6").Value = 1.03
$ws.Range("V6").Value = 1.16
$ws.Range("AD6").Value = 1000
$ws.Range("AI6").Value = 1000
$ws.Range("AJ6").Value = 1000
$ws.Range("F7").Value = 1.07
$ws.Range("K7").Value = 27
$ws.Range("N7").Value = 3.45
$ws.Range("O7").Value = 1.11
$ws.Range("Q7").Value = 1.35
$ws.Range("R7").Value = 1.92
$ws.Range("S7").Value = 1.77
$ws.Range("T7").Value = 2.98
$ws.Range("U7").Value = 1.4
$ws.Range("W7").Value = 1.99
$ws.Range("Y7").Value = 1000
$ws.Range("AF7").Value = 9.4
$ws.Range("AN7").Value = 2.5
$ws.Range("AB8").Value = 7.4
$ws.Range("AC8").Value = 10
$ws.Range("AD8").Value = 32
$ws.Range("AL8").Value = 55
$ws.Range("AH9").Value = 1000
$ws.Range("AL9").Value = 1000
$ws.Range("F10").Value = 1.25
$ws.Range("G10").Value = 1.3
$ws.Range("H10").Value = 1.09
$ws.Range("I10").Value = 22
$ws.Range("J10").Value = 4.7
$ws.Range("W10").Value = 3.9
$ws.Range("X10").Value = 22
$ws.Range("AF10").Value = 8
$ws.Range("AJ10").Value = 10.5
$ws.Range("G11").Value = 2.68
$ws.Range("H11").Value = 3.05
$ws.Range("I11").Value = 3.85
$ws.Range("J11").Value = 2.82
$ws.Range("N11").Value = 2.92
$ws.Range("O11").Value = 1.41
$ws.Range("P11").Value = 1.6
$ws.Range("R11").Value = 1.23
$ws.Range("U11").Value = 1.87
$ws.Range("W11").Value = 1.59
$ws.Range("X11").Value = 1000
$ws.Range("Y11").Value = 1000
$ws.Range("Z11").Value = 24
$ws.Range("AB11").Value = 10.5
$ws.Range("AC11").Value = 1000
$ws.Range("AM11").Value = 180
